$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A50:I50").Interior.Color = 12695295
